# Extended unified config process to EAs and OPs
#
# - Rename the "DropDownLists" sheet to "DropdownLists"
# - Make the (renamed) DropdownLists sheet the active/selected tab
#   (ExtendedAttributes therefore loses tabSelected)
# - Update the selection on DropdownLists to cell B35
#   (selection on ExtendedAttributes, E23, is left as-is)

$wb = $excel.ActiveWorkbook

$eaSheet  = $wb.Worksheets.Item("ExtendedAttributes")
$ddlSheet = $wb.Worksheets.Item("DropDownLists")

# Rename "DropDownLists" -> "DropdownLists"
$ddlSheet.Name = "DropdownLists"

# Switch the active sheet/tab to DropdownLists and move its selection to B35
$ddlSheet.Activate()
$ddlSheet.Range("B35").Select()
